$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / first tab)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5974
$ws1.Range("F5").Value = 82
$ws1.Range("F6").Value = 110
$ws1.Range("F9").Value = 551
$ws1.Range("F10").Value = 28

# Sheet "全部类型" (sheet4 / fourth tab)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5974
$ws4.Range("F6").Value = 82
$ws4.Range("F7").Value = 110
$ws4.Range("F11").Value = 551
$ws4.Range("F12").Value = 28
